$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, centered) from AC1 to the new header cells
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill in team record data for rows 2-45
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = 77  # AD
    $ws.Cells.Item($row, 31).Value = 85  # AE
    $ws.Cells.Item($row, 32).Value = 0   # AF
}
